$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-missing value in C3 (was 11.2)
$ws.Range("C3").Value = ""

# Remove row 26 ("RM 232") entirely - remaining rows shift up
$ws.Rows.Item(26).Delete()

# Remove the row that is now "SC 92" (originally row 28, now row 27 after the
# previous deletion) - remaining rows shift up again
$ws.Rows.Item(27).Delete()

# Update the missing-value pattern on the shifted rows to match the new
# "after" state of the dataset.

# Row 26 is now "SC 5": fill in column B (was blank)
$ws.Range("B26").Value = -20.2

# Row 27 is now "SC 101": column B becomes the new missing value
$ws.Range("B27").Value = ""

# Row 33 is now "SC 232": fill in columns B and C (were blank)
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
